$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the latest scraped cryptos data (prices + 1h volume deltas).
# Price cells that look numeric need an explicit text format so Excel
# keeps storing them as text (matching the source "35.614.12"-style strings)
# instead of silently coercing them into real numbers.

$ws.Range("D2").Value = "35.597.52"
$ws.Range("E2").Value = "  -2.94%  "

$ws.Range("D3").Value = "1.981.31"
$ws.Range("E3").Value = "  -3.87%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.04"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  -4.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.53"
$ws.Range("E7").Value = "  +7.03%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.81"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.360"
$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0736"
$ws.Range("E11").Value = "  -2.10%  "

$ws.Range("E12").Value = "  -2.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.939"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.61"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").Value = "2.272.24"
$ws.Range("E15").Value = "  -3.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.31"
$ws.Range("E16").Value = "  -3.13%  "

$ws.Range("D17").Value = "1.973.27"
$ws.Range("E17").Value = "  -4.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.98"
$ws.Range("E18").Value = "  +5.25%  "

$ws.Range("D19").Value = "35.538.88"
$ws.Range("E19").Value = "  -2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.33"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("D21").Value = "0.0₃0845"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.22"
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.45"
$ws.Range("E23").Value = "  -2.52%  "

$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.63"
$ws.Range("E25").Value = "  +22.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  -3.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.76"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.14"
$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("E29").Value = "  -5.04%  "

$ws.Range("E30").Value = "  -2.86%  "

$ws.Range("E31").Value = "  -4.63%  "

$ws.Range("E32").Value = "  -6.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0982"
$ws.Range("E33").Value = "  +18.74%  "

$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +9.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.34"
$ws.Range("E36").Value = "  -3.52%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("E38").Value = "  -3.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.32"
$ws.Range("E39").Value = "  +8.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.23"
$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.50"
$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.77"
$ws.Range("E44").Value = "  +2.77%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.09"
$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.26"
$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0896"
$ws.Range("E47").Value = "  -1.74%  "

$ws.Range("D48").Value = "1.372.49"
$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  -0.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.11"
$ws.Range("E50").Value = "  +2.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("E51").Value = "  -0.29%  "
